$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 1713.6666
$ws.Cells.Item(29, 9).Value = 86.5
$ws.Cells.Item(29, 10).Value = 2178.5715
$ws.Cells.Item(29, 11).Value = 259.5
$ws.Cells.Item(29, 12).Value = 6535.7145
$ws.Cells.Item(29, 13).Value = 21.5
$ws.Cells.Item(29, 14).Value = -7097.7145
$ws.Cells.Item(38, 8).Value = 688.7273
$ws.Cells.Item(38, 9).Value = 64
$ws.Cells.Item(38, 11).Value = 192
$ws.Cells.Item(38, 13).Value = 180
$ws.Cells.Item(58, 8).Value = 232.66667
$ws.Cells.Item(58, 9).Value = 232.66667
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 698.00001
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = -548.00001
$ws.Cells.Item(58, 14).ClearContents()
$ws.Cells.Item(87, 8).Value = 40050
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 40050
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 40050
$ws.Cells.Item(87, 13).ClearContents()
$ws.Cells.Item(87, 14).Value = -42546
$ws.Cells.Item(90, 8).Value = 40050
$ws.Cells.Item(90, 9).Value = 0
$ws.Cells.Item(90, 10).Value = 40050
$ws.Cells.Item(90, 11).Value = 0
$ws.Cells.Item(90, 12).Value = 120150
$ws.Cells.Item(90, 13).ClearContents()
$ws.Cells.Item(90, 14).Value = -132630
$ws.Cells.Item(106, 8).Value = 3288.4375
$ws.Cells.Item(106, 9).Value = 3288.4375
$ws.Cells.Item(106, 11).Value = 3288.4375
$ws.Cells.Item(106, 13).Value = -2657.4375
$ws.Cells.Item(129, 8).Value = 2096.2727
$ws.Cells.Item(129, 9).Value = 508.5
$ws.Cells.Item(129, 10).Value = 2315.276
$ws.Cells.Item(129, 11).Value = 1525.5
$ws.Cells.Item(129, 12).Value = 6945.828
$ws.Cells.Item(129, 13).Value = 3474.5
$ws.Cells.Item(129, 14).Value = -16945.828
$ws.Cells.Item(137, 8).Value = 1494.1
$ws.Cells.Item(137, 9).Value = 1244
$ws.Cells.Item(137, 10).Value = 1770.5264
$ws.Cells.Item(137, 11).Value = 3732
$ws.Cells.Item(137, 12).Value = 5311.5792
$ws.Cells.Item(137, 13).Value = -1182
$ws.Cells.Item(137, 14).Value = -10411.5792
# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 12724.174
$ws.Cells.Item(32, 9).Value = 3355.907
$ws.Cells.Item(32, 10).Value = 28217.846
$ws.Cells.Item(32, 11).Value = 3355.907
$ws.Cells.Item(32, 12).Value = 28217.846
$ws.Cells.Item(32, 13).Value = -3068.907
$ws.Cells.Item(32, 14).Value = -28791.846
$ws.Cells.Item(61, 8).Value = 3156
$ws.Cells.Item(61, 9).Value = 2399.4546
$ws.Cells.Item(61, 10).Value = 4344.857
$ws.Cells.Item(61, 11).Value = 2399.4546
$ws.Cells.Item(61, 12).Value = 4344.857
$ws.Cells.Item(61, 13).Value = -2187.4546
$ws.Cells.Item(61, 14).Value = -4768.857
$ws.Cells.Item(74, 8).Value = 15153689
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 15153689
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 15153689
$ws.Cells.Item(74, 13).ClearContents()
$ws.Cells.Item(74, 14).Value = -15155437
$ws.Cells.Item(77, 8).Value = 15153689
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 15153689
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 75768445
$ws.Cells.Item(77, 13).ClearContents()
$ws.Cells.Item(77, 14).Value = -75777181
$ws.Cells.Item(110, 8).Value = 5602.037
$ws.Cells.Item(110, 9).Value = 6493.095
$ws.Cells.Item(110, 10).Value = 2483.3333
$ws.Cells.Item(110, 11).Value = 6493.095
$ws.Cells.Item(110, 12).Value = 2483.3333
$ws.Cells.Item(110, 13).Value = -4448.095
$ws.Cells.Item(110, 14).Value = -6573.3333
$ws.Cells.Item(136, 8).Value = 3156
$ws.Cells.Item(136, 9).Value = 2399.4546
$ws.Cells.Item(136, 10).Value = 4344.857
$ws.Cells.Item(136, 11).Value = 7198.3638
$ws.Cells.Item(136, 12).Value = 13034.571
$ws.Cells.Item(136, 13).Value = -4648.3638
$ws.Cells.Item(136, 14).Value = -18134.571
# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 11125.917
$ws.Cells.Item(82, 9).Value = 6492
$ws.Cells.Item(82, 10).Value = 25027.666
$ws.Cells.Item(82, 11).Value = 6492
$ws.Cells.Item(82, 12).Value = 25027.666
$ws.Cells.Item(82, 13).Value = -6109
$ws.Cells.Item(82, 14).Value = -25793.666
$ws.Cells.Item(85, 8).Value = 11125.917
$ws.Cells.Item(85, 9).Value = 6492
$ws.Cells.Item(85, 10).Value = 25027.666
$ws.Cells.Item(85, 11).Value = 6492
$ws.Cells.Item(85, 12).Value = 25027.666
$ws.Cells.Item(85, 13).Value = -5166
$ws.Cells.Item(85, 14).Value = -27679.666
# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6898159
$ws.Cells.Item(31, 9).Value = 11111936
$ws.Cells.Item(31, 10).Value = 5001959.5
$ws.Cells.Item(31, 11).Value = 11111936
$ws.Cells.Item(31, 12).Value = 5001959.5
$ws.Cells.Item(31, 13).Value = -11111641
$ws.Cells.Item(31, 14).Value = -5002549.5
$ws.Cells.Item(34, 8).Value = 6898159
$ws.Cells.Item(34, 9).Value = 11111936
$ws.Cells.Item(34, 10).Value = 5001959.5
$ws.Cells.Item(34, 11).Value = 11111936
$ws.Cells.Item(34, 12).Value = 5001959.5
$ws.Cells.Item(34, 13).Value = -11111734
$ws.Cells.Item(34, 14).Value = -5002363.5
$ws.Cells.Item(36, 8).Value = 21666.334
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 10).Value = 21666.334
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 12).Value = 21666.334
$ws.Cells.Item(36, 13).ClearContents()
$ws.Cells.Item(36, 14).Value = -22442.334
$ws.Cells.Item(40, 8).Value = 21666.334
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 21666.334
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 21666.334
$ws.Cells.Item(40, 13).ClearContents()
$ws.Cells.Item(40, 14).Value = -21986.334
$ws.Cells.Item(42, 8).Value = 10000
$ws.Cells.Item(42, 10).Value = 10000
$ws.Cells.Item(42, 12).Value = 10000
$ws.Cells.Item(42, 14).Value = -11186
$ws.Cells.Item(50, 8).Value = 19999.166
$ws.Cells.Item(50, 10).Value = 19999.166
$ws.Cells.Item(50, 12).Value = 19999.166
$ws.Cells.Item(50, 14).Value = -21249.166
$ws.Cells.Item(58, 8).Value = 1373
$ws.Cells.Item(58, 9).Value = 815.63635
$ws.Cells.Item(58, 10).Value = 2054.2222
$ws.Cells.Item(58, 11).Value = 815.63635
$ws.Cells.Item(58, 12).Value = 2054.2222
$ws.Cells.Item(58, 13).Value = -612.63635
$ws.Cells.Item(58, 14).Value = -2460.2222
$ws.Cells.Item(136, 8).Value = 1373
$ws.Cells.Item(136, 9).Value = 815.63635
$ws.Cells.Item(136, 10).Value = 2054.2222
$ws.Cells.Item(136, 11).Value = 2446.90905
$ws.Cells.Item(136, 12).Value = 6162.6666
$ws.Cells.Item(136, 13).Value = 103.0909499999998
$ws.Cells.Item(136, 14).Value = -11262.6666
$ws.Cells.Item(138, 8).Value = 30000
$ws.Cells.Item(138, 10).Value = 30000
$ws.Cells.Item(138, 12).Value = 30000
$ws.Cells.Item(138, 14).Value = -40280
# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1116.9846
$ws.Cells.Item(68, 10).Value = 1247.7347
$ws.Cells.Item(68, 12).Value = 3743.2041
$ws.Cells.Item(68, 14).Value = -5365.2041
$ws.Cells.Item(71, 8).Value = 1116.9846
$ws.Cells.Item(71, 10).Value = 1247.7347
$ws.Cells.Item(71, 12).Value = 11229.6123
$ws.Cells.Item(71, 14).Value = -19341.6123
$ws.Cells.Item(131, 8).Value = 1026450.06
$ws.Cells.Item(131, 9).Value = 2469704.8
$ws.Cells.Item(131, 10).Value = 979.6053000000001
$ws.Cells.Item(131, 11).Value = 7409114.399999999
$ws.Cells.Item(131, 12).Value = 2938.8159
$ws.Cells.Item(131, 13).Value = -7404074.399999999
$ws.Cells.Item(131, 14).Value = -13018.8159
# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3030966.5
$ws.Cells.Item(22, 9).Value = 6667166.5
$ws.Cells.Item(22, 10).Value = 800
$ws.Cells.Item(22, 11).Value = 6667166.5
$ws.Cells.Item(22, 12).Value = 800
$ws.Cells.Item(22, 13).Value = -6666871.5
$ws.Cells.Item(22, 14).Value = -1390
$ws.Cells.Item(27, 8).Value = 3030966.5
$ws.Cells.Item(27, 9).Value = 6667166.5
$ws.Cells.Item(27, 10).Value = 800
$ws.Cells.Item(27, 11).Value = 6667166.5
$ws.Cells.Item(27, 12).Value = 800
$ws.Cells.Item(27, 13).Value = -6667059.5
$ws.Cells.Item(27, 14).Value = -1014
$ws.Cells.Item(61, 8).Value = 2365.3845
$ws.Cells.Item(61, 9).Value = 1975
$ws.Cells.Item(61, 10).Value = 3666.6667
$ws.Cells.Item(61, 11).Value = 1975
$ws.Cells.Item(61, 12).Value = 3666.6667
$ws.Cells.Item(61, 13).Value = -1773
$ws.Cells.Item(61, 14).Value = -4070.6667
$ws.Cells.Item(113, 8).Value = 2365.3845
$ws.Cells.Item(113, 9).Value = 1975
$ws.Cells.Item(113, 10).Value = 3666.6667
$ws.Cells.Item(113, 11).Value = 1975
$ws.Cells.Item(113, 12).Value = 3666.6667
$ws.Cells.Item(113, 13).Value = 195
$ws.Cells.Item(113, 14).Value = -8006.6667
